$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""

# Row 17
$ws.Range("H17").Value = 387347.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 387347.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1162043.4
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -1162379.4

# Row 40
$ws.Range("H40").Value = 1465
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""

# Row 137
$ws.Range("H137").Value = 2449
$ws.Range("I137").Value = 2298.8
$ws.Range("K137").Value = 6896.400000000001
$ws.Range("M137").Value = -4346.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2835.818
$ws.Range("I63").Value = 2274.25
$ws.Range("J63").Value = 4333.3335
$ws.Range("K63").Value = 2274.25
$ws.Range("L63").Value = 4333.3335
$ws.Range("M63").Value = -1588.25
$ws.Range("N63").Value = -5705.3335

# Row 66
$ws.Range("H66").Value = 2835.818
$ws.Range("I66").Value = 2274.25
$ws.Range("J66").Value = 4333.3335
$ws.Range("K66").Value = 11371.25
$ws.Range("L66").Value = 21666.6675
$ws.Range("M66").Value = -7939.25
$ws.Range("N66").Value = -28530.6675

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = ""

# Row 102
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622

# Row 132
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2613.6428
$ws.Range("I86").Value = 2063.8333
$ws.Range("J86").Value = 3026
$ws.Range("K86").Value = 2063.8333
$ws.Range("L86").Value = 3026
$ws.Range("M86").Value = -940.8332999999998
$ws.Range("N86").Value = -5272

# Row 89
$ws.Range("H89").Value = 2613.6428
$ws.Range("I89").Value = 2063.8333
$ws.Range("J89").Value = 3026
$ws.Range("K89").Value = 10319.1665
$ws.Range("L89").Value = 15130
$ws.Range("M89").Value = -4703.166499999999
$ws.Range("N89").Value = -26362

# Row 92
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992

# Row 94
$ws.Range("H94").Value = 2440.4443
$ws.Range("J94").Value = 1000
$ws.Range("L94").Value = 1000
$ws.Range("N94").Value = -1902

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = ""

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 350000
$ws.Range("I4").Value = 350000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 350000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -349888
$ws.Range("N4").Value = ""

# Row 31
$ws.Range("H31").Value = 6666.3335
$ws.Range("J31").Value = 10000
$ws.Range("L31").Value = 10000
$ws.Range("N31").Value = -10590

# Row 32
$ws.Range("H32").Value = 1997.6666
$ws.Range("I32").Value = 1997.6666
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1997.6666
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
$ws.Range("M32").Value = -1681.6666

# Row 34
$ws.Range("H34").Value = 6666.3335
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10404

# Row 93
$ws.Range("H93").Value = 5349.25
$ws.Range("I93").Value = 5349.25
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5349.25
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3477.25
$ws.Range("N93").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 4631533
$ws.Range("I122").Value = 4809631
$ws.Range("J122").Value = 995
$ws.Range("K122").Value = 14428893
$ws.Range("L122").Value = 2985
$ws.Range("M122").Value = -14426443
$ws.Range("N122").Value = -7885

# Row 132
$ws.Range("H132").Value = 4500
$ws.Range("J132").Value = 4500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3126.5715
$ws.Range("J7").Value = 2650
$ws.Range("L7").Value = 2650
$ws.Range("N7").Value = -2874

# Row 16
$ws.Range("H16").Value = 973
$ws.Range("I16").Value = 1146.5
$ws.Range("J16").Value = 799.5
$ws.Range("K16").Value = 1146.5
$ws.Range("L16").Value = 799.5
$ws.Range("M16").Value = -976.5
$ws.Range("N16").Value = -1139.5

# Row 122
$ws.Range("H122").Value = 3644.2856
$ws.Range("I122").Value = 3337
$ws.Range("J122").Value = 3874.75
$ws.Range("K122").Value = 10011
$ws.Range("L122").Value = 11624.25
$ws.Range("M122").Value = -7561
$ws.Range("N122").Value = -16524.25

# Row 126
$ws.Range("H126").Value = 3126.5715
$ws.Range("J126").Value = 2650
$ws.Range("L126").Value = 7950
$ws.Range("N126").Value = -12890

# Row 132
$ws.Range("H132").Value = 13874.875
$ws.Range("I132").Value = 6499.75
$ws.Range("J132").Value = 21250
$ws.Range("K132").Value = 19499.25
$ws.Range("L132").Value = 63750
$ws.Range("M132").Value = -16969.25
$ws.Range("N132").Value = -68810

# Row 136
$ws.Range("H136").Value = 3483.2
$ws.Range("I136").Value = 3483.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10449.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7899.599999999999
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1126.3334
$ws.Range("I96").Value = 1090
$ws.Range("J96").Value = 1144.5
$ws.Range("K96").Value = 1090
$ws.Range("L96").Value = 1144.5
$ws.Range("M96").Value = 283
$ws.Range("N96").Value = -3890.5

# Row 122
$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -9400

# Row 126
$ws.Range("H126").Value = 729.75
$ws.Range("I126").Value = 889.6667
$ws.Range("J126").Value = 250
$ws.Range("K126").Value = 2669.0001
$ws.Range("L126").Value = 750
$ws.Range("M126").Value = -199.0001000000002
$ws.Range("N126").Value = -5690

# Row 139
$ws.Range("H139").Value = 245000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
